$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 2
